# Rename disp_event IDs (column J) across all rows: "SD..." -> "SDTaq..."
# (keeps the "FT..." id_feces-family ids in column B/A untouched; only the
# seed-dispersal event id column is affected.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 10).End(-4162).Row  # xlUp = -4162

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 10)
    $old = $cell.Value2
    if ($old -ne $null -and $old.ToString().StartsWith("SD")) {
        $new = "SDTaq" + $old.ToString().Substring(2)
        $cell.Value = $new
    }
}

# Restore the active-cell selection to J1 (matches the saved sheetView state)
$ws.Range("J1").Select()
